$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New client row (row 4): Niles Grill, Banned ---
# Pre-format the purely-numeric-looking cells as Text so Excel stores them
# as strings (matching "6666"/"4039013212" style already used by the other
# id/phone columns) instead of auto-converting them to numbers.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"

$ws.Range("A4").Value = "6666"
$ws.Range("B4").Value = "Grill"
$ws.Range("C4").Value = "Niles"
$ws.Range("D4").Value = "4039013212"
$ws.Range("E4").Value = "nile.grill@gmail.com"
$ws.Range("F4").Value = "Banned"

# --- New equipment / rentals block (rows 6-8, columns H:L) ---
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"

$ws.Range("H6").Value = "12345"
$ws.Range("I6").Value = "Pickaxe"
$ws.Range("J6").Value = "Tool"
$ws.Range("K6").Value = "Mining"
$ws.Range("L6").Value = 12

$ws.Range("H7").Value = "54321"
$ws.Range("I7").Value = "Shovel"
$ws.Range("J7").Value = "Tool"
$ws.Range("K7").Value = "Dig Stuff"
$ws.Range("L7").Value = 12

$ws.Range("H8").Value = "15243"
$ws.Range("I8").Value = "Helmet"
$ws.Range("J8").Value = "Safety"
$ws.Range("K8").Value = "Protects Head"
$ws.Range("L8").Value = 5

# Widen the new "Description" column (K) like the other text columns.
$ws.Columns.Item(11).ColumnWidth = 24.166666666666668

# Match the updated view: zoomed out a bit and selection parked on E7.
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("E7").Select() | Out-Null
